# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) updates: row -> new value
$exhibitionUpdates = @{
    2  = 6842
    4  = 438
    5  = 72
    7  = 547
    8  = 111
    11 = 6
    12 = 33
    13 = 182
    14 = 427
    15 = 5
    16 = 1630
    17 = 28
    18 = 3458
    20 = 236
    21 = 15
    22 = 2093
    23 = 178
    24 = 6
    26 = 4
    28 = 12
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (sheet4) updates: row -> new value
$allTypesUpdates = @{
    2  = 6842
    4  = 438
    5  = 72
    8  = 547
    9  = 111
    12 = 6
    13 = 33
    14 = 182
    15 = 427
    16 = 5
    17 = 1630
    18 = 28
    19 = 3458
    21 = 236
    22 = 15
    23 = 2093
    24 = 178
    25 = 6
    27 = 4
    29 = 12
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
